$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "def"
$ws.Range("B2").Value = "KEYWORD"
$ws.Range("C2").Value = 1

$ws.Range("A3").Value = "int"
$ws.Range("B3").Value = "KEYWORD"
$ws.Range("C3").Value = 1

$ws.Range("A4").Value = "gcd"
$ws.Range("B4").Value = "ID"
$ws.Range("C4").Value = 1

$ws.Range("A5").Value = "("
$ws.Range("B5").Value = "STR_BKT"
$ws.Range("C5").Value = 1

$ws.Range("A6").Value = "double"
$ws.Range("B6").Value = "KEYWORD"
$ws.Range("C6").Value = 1

$ws.Range("A7").Value = "a"
$ws.Range("B7").Value = "ID"
$ws.Range("C7").Value = 1

$ws.Range("A8").Value = ","
$ws.Range("B8").Value = "COM"
$ws.Range("C8").Value = 1

$ws.Range("A9").Value = "int"
$ws.Range("B9").Value = "KEYWORD"
$ws.Range("C9").Value = 1

$ws.Range("A10").Value = "b"
$ws.Range("B10").Value = "ID"
$ws.Range("C10").Value = 1

$ws.Range("A11").Value = ")"
$ws.Range("B11").Value = "END_BKT"
$ws.Range("C11").Value = 1

$ws.Range("A12").Value = "if"
$ws.Range("B12").Value = "KEYWORD"
$ws.Range("C12").Value = 2

$ws.Range("A13").Value = "("
$ws.Range("B13").Value = "STR_BKT"
$ws.Range("C13").Value = 2

$ws.Range("A14").Value = "a"
$ws.Range("B14").Value = "ID"
$ws.Range("C14").Value = 2

$ws.Range("A15").Formula = "=="
$ws.Range("B15").Value = "EV"
$ws.Range("C15").Value = 2

$ws.Range("A16").Value = "b"
$ws.Range("B16").Value = "ID"
$ws.Range("C16").Value = 2

$ws.Range("A17").Value = ")"
$ws.Range("B17").Value = "END_BKT"
$ws.Range("C17").Value = 2

$ws.Range("A18").Value = "then"
$ws.Range("B18").Value = "KEYWORD"
$ws.Range("C18").Value = 2

$ws.Range("A19").Value = "return"
$ws.Range("B19").Value = "KEYWORD"
$ws.Range("C19").Value = 2

$ws.Range("A20").Value = "("
$ws.Range("B20").Value = "STR_BKT"
$ws.Range("C20").Value = 2

$ws.Range("A21").Value = "a"
$ws.Range("B21").Value = "ID"
$ws.Range("C21").Value = 2

$ws.Range("A22").Value = ")"
$ws.Range("B22").Value = "END_BKT"
$ws.Range("C22").Value = 2

$ws.Range("A23").Value = "fi"
$ws.Range("B23").Value = "KEYWORD"
$ws.Range("C23").Value = 2

$ws.Range("A24").Value = ";"
$ws.Range("B24").Value = "SEMI"
$ws.Range("C24").Value = 2

$ws.Range("A25").Value = "if"
$ws.Range("B25").Value = "KEYWORD"
$ws.Range("C25").Value = 3

$ws.Range("A26").Value = "("
$ws.Range("B26").Value = "STR_BKT"
$ws.Range("C26").Value = 3

$ws.Range("A27").Value = "a"
$ws.Range("B27").Value = "ID"
$ws.Range("C27").Value = 3

$ws.Range("A28").Value = ">"
$ws.Range("B28").Value = "GT"
$ws.Range("C28").Value = 3

$ws.Range("A29").Value = "b"
$ws.Range("B29").Value = "ID"
$ws.Range("C29").Value = 3

$ws.Range("A30").Value = ")"
$ws.Range("B30").Value = "END_BKT"
$ws.Range("C30").Value = 3

$ws.Range("A31").Value = "then"
$ws.Range("B31").Value = "KEYWORD"
$ws.Range("C31").Value = 3

$ws.Range("A32").Value = "return"
$ws.Range("B32").Value = "KEYWORD"
$ws.Range("C32").Value = 3

$ws.Range("A33").Value = "("
$ws.Range("B33").Value = "STR_BKT"
$ws.Range("C33").Value = 3

$ws.Range("A34").Value = "gcd"
$ws.Range("B34").Value = "ID"
$ws.Range("C34").Value = 3

$ws.Range("A35").Value = "("
$ws.Range("B35").Value = "STR_BKT"
$ws.Range("C35").Value = 3

$ws.Range("A36").Value = "a"
$ws.Range("B36").Value = "ID"
$ws.Range("C36").Value = 3

$ws.Range("A37").Value = "-"
$ws.Range("B37").Value = "SUB"
$ws.Range("C37").Value = 3

$ws.Range("A38").Value = "b"
$ws.Range("B38").Value = "ID"
$ws.Range("C38").Value = 3

$ws.Range("A39").Value = ","
$ws.Range("B39").Value = "COM"
$ws.Range("C39").Value = 3

$ws.Range("A40").Value = "b"
$ws.Range("B40").Value = "ID"
$ws.Range("C40").Value = 3

$ws.Range("A41").Value = ")"
$ws.Range("B41").Value = "END_BKT"
$ws.Range("C41").Value = 3

$ws.Range("A42").Value = ")"
$ws.Range("B42").Value = "END_BKT"
$ws.Range("C42").Value = 3

$ws.Range("A43").Value = "else"
$ws.Range("B43").Value = "KEYWORD"
$ws.Range("C43").Value = 4

$ws.Range("A44").Value = "return"
$ws.Range("B44").Value = "KEYWORD"
$ws.Range("C44").Value = 4

$ws.Range("A45").Value = "("
$ws.Range("B45").Value = "STR_BKT"
$ws.Range("C45").Value = 4

$ws.Range("A46").Value = "gcd"
$ws.Range("B46").Value = "ID"
$ws.Range("C46").Value = 4

$ws.Range("A47").Value = "("
$ws.Range("B47").Value = "STR_BKT"
$ws.Range("C47").Value = 4

$ws.Range("A48").Value = "a"
$ws.Range("B48").Value = "ID"
$ws.Range("C48").Value = 4

$ws.Range("A49").Value = ","
$ws.Range("B49").Value = "COM"
$ws.Range("C49").Value = 4

$ws.Range("A50").Value = "b"
$ws.Range("B50").Value = "ID"
$ws.Range("C50").Value = 4

$ws.Range("A51").Value = "-"
$ws.Range("B51").Value = "SUB"
$ws.Range("C51").Value = 4

$ws.Range("A52").Value = "a"
$ws.Range("B52").Value = "ID"
$ws.Range("C52").Value = 4

$ws.Range("A53").Value = ")"
$ws.Range("B53").Value = "END_BKT"
$ws.Range("C53").Value = 4

$ws.Range("A54").Value = ")"
$ws.Range("B54").Value = "END_BKT"
$ws.Range("C54").Value = 4

$ws.Range("A55").Value = "fi"
$ws.Range("B55").Value = "KEYWORD"
$ws.Range("C55").Value = 4

$ws.Range("A56").Value = ";"
$ws.Range("B56").Value = "SEMI"
$ws.Range("C56").Value = 4

$ws.Range("A57").Value = "fed"
$ws.Range("B57").Value = "KEYWORD"
$ws.Range("C57").Value = 5

$ws.Range("A58").Value = ";"
$ws.Range("B58").Value = "SEMI"
$ws.Range("C58").Value = 5

$ws.Range("A59").Value = "print"
$ws.Range("B59").Value = "KEYWORD"
$ws.Range("C59").Value = 6

$ws.Range("A60").Value = "gcd"
$ws.Range("B60").Value = "ID"
$ws.Range("C60").Value = 6

$ws.Range("A61").Value = "("
$ws.Range("B61").Value = "STR_BKT"
$ws.Range("C61").Value = 6

$ws.Range("A62").NumberFormat = "@"
$ws.Range("A62").Value = "21"
$ws.Range("A62").Style = "Normal"
$ws.Range("B62").Value = "INT"
$ws.Range("C62").Value = 6

$ws.Range("A63").Value = ","
$ws.Range("B63").Value = "COM"
$ws.Range("C63").Value = 6

$ws.Range("A64").NumberFormat = "@"
$ws.Range("A64").Value = "15"
$ws.Range("A64").Style = "Normal"
$ws.Range("B64").Value = "INT"
$ws.Range("C64").Value = 6

$ws.Range("A65").Value = ")"
$ws.Range("B65").Value = "END_BKT"
$ws.Range("C65").Value = 6

$ws.Range("A66").Value = ";"
$ws.Range("B66").Value = "SEMI"
$ws.Range("C66").Value = 6

$ws.Range("A67").Value = "print"
$ws.Range("B67").Value = "KEYWORD"
$ws.Range("C67").Value = 6

$ws.Range("A68").NumberFormat = "@"
$ws.Range("A68").Value = "45"
$ws.Range("A68").Style = "Normal"
$ws.Range("B68").Value = "INT"
$ws.Range("C68").Value = 6

$ws.Range("A69").Value = ";"
$ws.Range("B69").Value = "SEMI"
$ws.Range("C69").Value = 6

$ws.Range("A70").Value = "print"
$ws.Range("B70").Value = "KEYWORD"
$ws.Range("C70").Value = 6

$ws.Range("A71").NumberFormat = "@"
$ws.Range("A71").Value = "2"
$ws.Range("A71").Style = "Normal"
$ws.Range("B71").Value = "INT"
$ws.Range("C71").Value = 6

$ws.Range("A72").Value = "*"
$ws.Range("B72").Value = "MLT"
$ws.Range("C72").Value = 6

$ws.Range("A73").Value = "("
$ws.Range("B73").Value = "STR_BKT"
$ws.Range("C73").Value = 6

$ws.Range("A74").Value = "gcd"
$ws.Range("B74").Value = "ID"
$ws.Range("C74").Value = 6

$ws.Range("A75").Value = "("
$ws.Range("B75").Value = "STR_BKT"
$ws.Range("C75").Value = 6

$ws.Range("A76").NumberFormat = "@"
$ws.Range("A76").Value = "21"
$ws.Range("A76").Style = "Normal"
$ws.Range("B76").Value = "INT"
$ws.Range("C76").Value = 6

$ws.Range("A77").Value = ","
$ws.Range("B77").Value = "COM"
$ws.Range("C77").Value = 6

$ws.Range("A78").NumberFormat = "@"
$ws.Range("A78").Value = "28"
$ws.Range("A78").Style = "Normal"
$ws.Range("B78").Value = "INT"
$ws.Range("C78").Value = 6

$ws.Range("A79").Value = ")"
$ws.Range("B79").Value = "END_BKT"
$ws.Range("C79").Value = 6

$ws.Range("A80").Value = "+"
$ws.Range("B80").Value = "ADD"
$ws.Range("C80").Value = 6

$ws.Range("A81").NumberFormat = "@"
$ws.Range("A81").Value = "6"
$ws.Range("A81").Style = "Normal"
$ws.Range("B81").Value = "INT"
$ws.Range("C81").Value = 6

$ws.Range("A82").Value = ")"
$ws.Range("B82").Value = "END_BKT"
$ws.Range("C82").Value = 6

$ws.Range("A83").Value = "."
$ws.Range("B83").Value = "PER"
$ws.Range("C83").Value = 6

$ws.Range("A84:C98").Clear()
